# Add two new columns, I (I0) and J (IF), to the worksheet.
# I column: header "I0", values = 1 for every data row (2..33)
# J column: header "IF", values = same as column H for every data row (2..33)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 33

# Headers (copy style from existing header cell H1 so the new headers match
# the bold/centered/bordered header style already used in row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: I = 1, J = same value as H
for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
